# Update the KPI computations across several instance-data sheets.
$wb = $excel.ActiveWorkbook

# ---- Productdata sheet ----
$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 14.52533333333333
$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 5.861066666666667
$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 3.816
$ws.Range("C5").Value = 0
$ws.Range("E5").Value = 1.59
$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 2.862
$ws.Range("C7").Value = 718
$ws.Range("E7").Value = 12.194
$ws.Range("C8").Value = 143
$ws.Range("E8").Value = 4.555066666666668
$ws.Range("C9").Value = 234
$ws.Range("E9").Value = 3.637333333333333

# ---- ForecastedAverageDemand sheet ----
$ws = $wb.Worksheets.Item("ForecastedAverageDemand")
$ws.Range("G2").Value = 242
$ws.Range("H2").Value = 44
$ws.Range("G3").Value = 199
$ws.Range("H3").Value = 45
$ws.Range("I3").Value = 72
$ws.Range("G4").Value = 237
$ws.Range("H4").Value = 46
$ws.Range("G5").Value = 232
$ws.Range("H5").Value = 45
$ws.Range("G6").Value = 227
$ws.Range("H6").Value = 44
$ws.Range("I6").Value = 73
$ws.Range("G7").Value = 270
$ws.Range("I7").Value = 72

# ---- ForcastedStandardDeviation sheet ----
$ws = $wb.Worksheets.Item("ForcastedStandardDeviation")
$ws.Range("G2").Value = 6.049999999999999
$ws.Range("H2").Value = 1.1
$ws.Range("G3").Value = 9.452499999999997
$ws.Range("H3").Value = 2.137499999999999
$ws.Range("I3").Value = 3.419999999999999
$ws.Range("G4").Value = 16.05674999999999
$ws.Range("H4").Value = 3.116499999999999
$ws.Range("G5").Value = 19.9462
$ws.Range("H5").Value = 3.868875
$ws.Range("G6").Value = 23.2396925
$ws.Range("H6").Value = 4.50461
$ws.Range("I6").Value = 7.473557499999998
$ws.Range("G7").Value = 31.6277325
$ws.Range("I7").Value = 8.434061999999999

# ---- Capacity sheet ----
$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("B2").Value = 2793.333333333333
$ws.Range("B3").Value = 5908.333333333332
$ws.Range("B4").Value = 7950
$ws.Range("B5").Value = 7950
$ws.Range("B6").Value = 3975
$ws.Range("B7").Value = 11725
$ws.Range("B8").Value = 448.3333333333334
$ws.Range("B9").Value = 2200

# ---- ProcessingTime sheet ----
$ws = $wb.Worksheets.Item("ProcessingTime")
$ws.Range("B2").Value = 1
$ws.Range("C3").Value = 5
$ws.Range("D4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G7").Value = 5
$ws.Range("I9").Value = 3
